$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# ---------------------------------------------------------------
# Sheet "Sale 22-23": insert the new rows (structure first)
# ---------------------------------------------------------------

# New row 13 (additional payment entry between the row-12 subtotal and row-14 block)
$ws2.Rows.Item(13).Insert()

# Two new rows before the old header-row-24 (now at row 25 after the shift above)
$ws2.Rows.Item(25).Insert()
$ws2.Rows.Item(26).Insert()

# Restore formatting for the freshly inserted rows by pasting formats from
# sibling rows that carry the same visual pattern.
$ws2.Range("A12:F12").Copy() | Out-Null
$ws2.Range("A13:F13").PasteSpecial(-4122) | Out-Null

$ws2.Range("A23:F23").Copy() | Out-Null
$ws2.Range("A25:F25").PasteSpecial(-4122) | Out-Null
$ws2.Application.CutCopyMode = $false

# Drop the now-superseded subtotal formula in F12 (stays blank, keeps its style)
$ws2.Range("F12").ClearContents()

# Row 13 data
$ws2.Range("B13").Value = 45038
$ws2.Range("C13").Value = "b23-24MQ107"
$ws2.Range("D13").Value = "Putzmeister Concrete Machines Pvt Ltd"
$ws2.Range("E13").Value = 151972
$ws2.Range("F13").Formula = "=E5+E6+E7+E8+E9+E10+E11+E12+E13"

# Row 25 data (entry 8)
$ws2.Range("A25").Value = 8
$ws2.Range("B25").Value = 45034
$ws2.Range("C25").Value = "b23-24MQ106"
$ws2.Range("D25").Value = "Marcfremiot"
$ws2.Range("E25").Value = 21900.8
$ws2.Range("F25").Formula = "=E25"

# ---------------------------------------------------------------
# Sheet "Purchase 22-23": append the two new entries + the note row
# ---------------------------------------------------------------

# Row 21 (entry 7) - clone formatting from row 19, then fill in values
$ws1.Range("A19:F19").Copy() | Out-Null
$ws1.Range("A21:F21").PasteSpecial(-4122) | Out-Null
$ws1.Application.CutCopyMode = $false
$ws1.Rows.Item(21).RowHeight = 14.4

$ws1.Range("A21").Value = 7
$ws1.Range("B21").Value = 45029
$ws1.Range("C21").Value = "SLH/135"
$ws1.Range("D21").Value = "Shree Laxmi Lighting Hub"
$ws1.Range("E21").Value = 1540
$ws1.Range("F21").Formula = "=E21"

# Row 23 (entry 8) - clone formatting from row 19, then fill in values
$ws1.Range("A19:F19").Copy() | Out-Null
$ws1.Range("A23:F23").PasteSpecial(-4122) | Out-Null
$ws1.Application.CutCopyMode = $false
$ws1.Rows.Item(23).RowHeight = 14.4

$ws1.Range("A23").Value = 8
$ws1.Range("B23").Value = 45036
$ws1.Range("C23").Value = "Digi.2324/0090"
$ws1.Range("D23").Value = "Digiserve"
$ws1.Range("E23").Value = 93220
$ws1.Range("F23").Formula = "=E23"

# Row 29 - standalone note/total cell
$ws1.Range("D29").Formula = "=28933+885"

# ---------------------------------------------------------------
# Selections: Sale 22-23 -> F26, Purchase 22-23 -> E27 (stays the active tab)
# ---------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F26").Select()
$ws1.Activate()
$ws1.Range("E27").Select()
